$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the bill requester's identity / label fields
$ws.Range("A3").Value = "নাম: Dr. Sk. Md. Masudul Ahsan"
$ws.Range("A4").Value = "পদবী: অধ্যাপক"
$ws.Range("F5").Value = "বিভাগ :সিএসই"
$ws.Range("A32").Value = "কথায়:দুই লক্ষ বত্রিশ হাজার ছয়শত্লিশ টাকা মাত্র।"

# Fill in the quantities for each billed item (label-wise bill generation)
$ws.Range("G9").Value = 57
$ws.Range("G12").Value = 57
$ws.Range("G14").Value = 59
$ws.Range("G16").Value = 27
$ws.Range("G17").Value = 19.5
$ws.Range("G18").Value = 118
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 5
$ws.Range("G29").Value = 15

# Move the active selection to B5, like the author left it
$ws.Range("B5").Select()
